# "firebase finialization and clientui features"
#
# This edit adds a new "E" column of reviewer/assignee names to the
# ADMIN UI sheet (sheet1) and updates the saved view/selection state on
# both the ADMIN UI and CLIENT UI sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ADMIN UI")
$ws2 = $wb.Worksheets.Item("CLIENT UI")

# ---------------------------------------------------------------
# CLIENT UI sheet: only the view/selection changes (no cell data
# changes). Do this first so that ADMIN UI ends up as the final
# active/selected sheet & tab, matching the original file.
# ---------------------------------------------------------------
$ws2.Select()
$excel.ActiveWindow.ScrollRow    = 1
$excel.ActiveWindow.ScrollColumn = 2
$ws2.Range("E12").Select()

# ---------------------------------------------------------------
# ADMIN UI sheet: new column E values.
# Insert the new shared strings in this order so they land at the
# same shared-string-table indices as the target workbook:
#   66 = "thaniga", 67 = "christen", 68 = "Dhamu"
# ---------------------------------------------------------------
$ws1.Select()

$ws1.Range("E4").Value = "thaniga"

# E25 should carry the same fill/style as D25 (s="5") -- copy the
# formatting over first, then overwrite the value with the new text.
$ws1.Range("D25").Copy($ws1.Range("E25"))
$ws1.Range("E25").Value = "christen"

$ws1.Range("E3").Value = "Dhamu"
$ws1.Range("E6").Value = "Dhamu"

# Restore the view/selection state for ADMIN UI.
$excel.ActiveWindow.ScrollRow    = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("E5").Select()
